$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in values for the "17" row (row 20): G20, H20 already have borders/style,
# I20 and J20 are new plain cells.
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 5

# Update the active selection to match the post-edit state (K20).
$ws.Range("K20").Select()
